# Movement through movement map nearly done
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Update B8 (e007 Morning Briefing - Weather Roll): insert "B29" after the leading LineBreaks
$ws.Range("B8").Value = "<Bold>e007 Morning Briefing - Weather Roll</Bold> <InlineUIContainer><Button Content='r4.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`n<LineBreak/><LineBreak/>B29`nThe `n<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `n Table determines weather for today:  `n<InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer>`n<LineBreak/>"

# Update B33 (e032 No Combat): swap Image for Button and drop trailing "Continue with "
$ws.Range("B33").Value = "<Bold>e032 No Combat</Bold> `n<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nNo combat. Victory points added to the After Action Report `n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nto reflect area under US Control. "

# Update the view - scroll down a bit
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 30
